# The workbook has two sheets: "personal computer" and "server computer".
# Both originally share the string "BlockingStlQueue" for their D1 header
# cell. This edit renames the header on the "server computer" sheet to
# "BlockingStlQueueWithMutex" (the "personal computer" sheet keeps its
# original "BlockingStlQueue" header), widens column D on the
# "personal computer" sheet to fit the new, longer label used elsewhere in
# the workbook, and updates the remembered selection on that sheet.

$wb = $excel.ActiveWorkbook

# 1) Rename the D1 header on the "server computer" sheet.
$wsServer = $wb.Worksheets.Item("server computer")
$wsServer.Range("D1").Value = "BlockingStlQueueWithMutex"

# 2) On the "personal computer" sheet: widen column D to 17 characters
#    (Excel's ColumnWidth property is expressed in characters of the
#    Normal style's font, which is offset from the stored column width by
#    5/MaximumDigitWidth; MaximumDigitWidth is 7 for the default Calibri 11
#    font, so subtract 5/7 to land on a stored width of exactly 17).
$wsPersonal = $wb.Worksheets.Item("personal computer")
$wsPersonal.Columns.Item(4).ColumnWidth = 17 - 5/7

# 3) Update the saved selection on the "personal computer" sheet to G33.
$wsPersonal.Activate()
$null = $wsPersonal.Range("G33").Select()
